$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Update Runmode column (C) to "N" for all suite rows except row 4 (C Suite)
$ws.Range("C2").Value = "N"
$ws.Range("C3").Value = "N"
$ws.Range("C5").Value = "N"
$ws.Range("C6").Value = "N"
$ws.Range("C7").Value = "N"

# Update the active selection to C4
$ws.Activate()
$ws.Range("C4").Select()
